$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(1.445647641019636, 1.626987699542094, 18.71679738969934, 13.86384647080068, 1, 35.65327920106175)
    3 = @(0.04172184405617529, 0.3048912486333797, 0.7210945179870265, 2797.565817734744, 1, 2798.63352534542)
    4 = @(1.445647641019636, 1.626987699542094, 189.6080260415259, 0.5333859586016987, 1, 193.2140473406893)
    5 = @(1.445647641019636, 1.626987699542094, 0.1496068669990043, 2797.565817734744, 0, 2800.788059942304)
    6 = @(0.6545652718822623, 0.3048912486333797, 0.1496068669990043, 0.5333859586016987, 1, 1.642449346116345)
    7 = @(3.272327238179451, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 1, 5.582307763322248)
    8 = @(0.6545652718822623, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 1, 2.964545797025059)
    9 = @(0.1169995834814548, 0.3048912486333797, 3.223369029078222, 13.86384647080068, 1, 17.50910633199374)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
}
